# Update "programs" sheet: add a DAYS column (E) describing which days of
# the week each program meets on, then make "programs" the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("programs")

# Header (bold, matching the other header cells in row 1)
$ws.Range("E1").Value = "DAYS"
$ws.Range("E1").Font.Bold = $true

# Data rows - days-of-week codes (1=Mon ... 7=Sun), pipe separated.
# Seed each distinct value once first (in the same order the source
# workbook introduced them) so shared-string indices line up, then
# backfill the repeated cells.
$ws.Range("E2").Value = "1|3|5"
$ws.Range("E7").Value = "1|3"
$ws.Range("E11").Value = "6|7"
$ws.Range("E5").Value = "1|2|3|4|5"
$ws.Range("E15").Value = "5|6|7"

$ws.Range("E3").Value = "1|3|5"
$ws.Range("E4").Value = "1|3|5"
$ws.Range("E6").Value = "1|2|3|4|5"
$ws.Range("E8").Value = "1|3"
$ws.Range("E9").Value = "1|3"
$ws.Range("E10").Value = "1|3"
$ws.Range("E12").Value = "6|7"
$ws.Range("E13").Value = "6|7"
$ws.Range("E14").Value = "6|7"
$ws.Range("E16").Value = "5|6|7"

# New column E width
$ws.Range("E1").ColumnWidth = 9.1

# Make "programs" the active sheet/tab, with E9 selected
$ws.Activate()
$ws.Range("E9").Select()
